# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 18:22"

# Row 4 (Estados Unidos) - updated counters
$ws.Range("B4").Value = 291545
$ws.Range("C4").Value = 14384
$ws.Range("E4").Value = 269326
$ws.Range("G4").Value = 447
$ws.Range("H4").Value = 7851

# Row 6 (Alemania) - updated counters
$ws.Range("B6").Value = 124632
$ws.Range("C6").Value = 4805
$ws.Range("D6").Value = 20996
$ws.Range("E6").Value = 88274
$ws.Range("F6").Value = 3994
$ws.Range("G6").Value = 681
$ws.Range("H6").Value = 15362

# Row 16 - updated counters
$ws.Range("B16").Value = 12924
$ws.Range("C16").Value = 549
$ws.Range("E16").Value = 10388

# Row 20 - updated counters
$ws.Range("B20").Value = 9391
$ws.Range("C20").Value = 197
$ws.Range("E20").Value = 8888
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 376

# Row 28 - updated counters
$ws.Range("D28").Value = 528
$ws.Range("E28").Value = 3606
$ws.Range("F28").Value = 38

# Row 31 - updated counters
$ws.Range("E31").Value = 3138
$ws.Range("G31").Value = 13
$ws.Range("H31").Value = 146

# Row 37 - updated counters
$ws.Range("B37").Value = 2748
$ws.Range("C37").Value = 62
$ws.Range("E37").Value = 2577

# Row 41 - updated counters
$ws.Range("D41").Value = 674
$ws.Range("E41").Value = 1373

# Argelia overtakes Singapur in the ranking (sorted descending by Casos totales):
# Row 54 becomes Argelia with newly updated data
$ws.Range("A54").Value = "Argelia"
$ws.Range("B54").Value = 1251
$ws.Range("C54").Value = 80
$ws.Range("D54").Value = 90
$ws.Range("E54").Value = 1031
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 25
$ws.Range("H54").Value = 130

# Row 55 becomes Singapur, keeping its previous data
$ws.Range("A55").Value = "Singapur"
$ws.Range("B55").Value = 1189
$ws.Range("C55").Value = 75
$ws.Range("D55").Value = 297
$ws.Range("E55").Value = 886
$ws.Range("F55").Value = 24
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 6

# Row 94 - updated counters
$ws.Range("F94").Value = 11

# Row 110 - updated counters
$ws.Range("D110").Value = 31
$ws.Range("E110").Value = 125
